$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3579
$ws.Range("F5").Value = 8349
$ws.Range("F8").Value = 2243
$ws.Range("F10").Value = 106
$ws.Range("F11").Value = 78
$ws.Range("F12").Value = 663
$ws.Range("F14").Value = 7515
$ws.Range("F15").Value = 449
$ws.Range("F16").Value = 7734
$ws.Range("F18").Value = 57930
$ws.Range("F19").Value = 57930
$ws.Range("F20").Value = 4847
$ws.Range("F21").Value = 1065
$ws.Range("F22").Value = 962
$ws.Range("F23").Value = 515
$ws.Range("F25").Value = 933
$ws.Range("F28").Value = 5321
$ws.Range("F30").Value = 123
$ws.Range("F31").Value = 53
$ws.Range("F32").Value = 923
$ws.Range("F33").Value = 1413
$ws.Range("F34").Value = 2003
$ws.Range("F35").Value = 22
$ws.Range("F36").Value = 190
$ws.Range("F37").Value = 234
$ws.Range("F38").Value = 1095
$ws.Range("F39").Value = 4
$ws.Range("F40").Value = 733
$ws.Range("F41").Value = 45
$ws.Range("F42").Value = 788
$ws.Range("F43").Value = 278
$ws.Range("F44").Value = 263
$ws.Range("F47").Value = 210

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 38
$ws.Range("F6").Value = 158
$ws.Range("G6").Value = 280
$ws.Range("F8").Value = 52
$ws.Range("F9").Value = 7669
$ws.Range("F10").Value = 127
$ws.Range("F14").Value = 8
$ws.Range("F20").Value = 24
$ws.Range("F23").Value = 36
$ws.Range("F27").Value = 132
$ws.Range("F33").Value = 8
$ws.Range("F38").Value = 51
$ws.Range("G41").Value = 880
$ws.Range("F45").Value = 41
$ws.Range("F47").Value = 68
$ws.Range("F48").Value = 283

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2384
$ws.Range("F5").Value = 1622
$ws.Range("F8").Value = 2429
$ws.Range("F9").Value = 9452
$ws.Range("F10").Value = 1774
$ws.Range("F11").Value = 183
$ws.Range("F15").Value = 280
$ws.Range("F16").Value = 2427
$ws.Range("F17").Value = 172
$ws.Range("F18").Value = 69
$ws.Range("F19").Value = 536

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 3579
$ws.Range("F4").Value = 2384
$ws.Range("F5").Value = 8349
$ws.Range("F6").Value = 1774
$ws.Range("F7").Value = 183
$ws.Range("F8").Value = 280
$ws.Range("F9").Value = 172
$ws.Range("F10").Value = 78
$ws.Range("F11").Value = 663
$ws.Range("F12").Value = 7734
$ws.Range("F13").Value = 57930
$ws.Range("F15").Value = 38
$ws.Range("F16").Value = 4848
$ws.Range("F17").Value = 962
$ws.Range("F18").Value = 515
$ws.Range("F19").Value = 933
$ws.Range("F21").Value = 158
$ws.Range("G21").Value = 280
$ws.Range("F22").Value = 123
$ws.Range("F23").Value = 923
$ws.Range("F24").Value = 1413
$ws.Range("F25").Value = 2003
$ws.Range("F26").Value = 127
$ws.Range("F27").Value = 536
$ws.Range("F29").Value = 8
$ws.Range("F32").Value = 234
$ws.Range("F33").Value = 24
$ws.Range("F34").Value = 45
$ws.Range("F35").Value = 788
$ws.Range("F36").Value = 278
$ws.Range("F40").Value = 210
$ws.Range("F41").Value = 8
$ws.Range("F44").Value = 51
$ws.Range("F49").Value = 41
$ws.Range("F50").Value = 68
